$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input")

# Update mobile number value and reset style to default (General)
$ws.Range("B5").Style = "Normal"
$ws.Range("B5").Value = 9898123456

# Update selection / active cell to E9
$ws.Range("E9").Select()
